$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.958.33'
$ws.Range("E2").Value = '  -0.15%  '

$ws.Range("D3").Value = '2.590.45'
$ws.Range("E3").Value = '  +1.48%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.86'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.43%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.15'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.36%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.596'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +2.15%  '

$ws.Range("E9").Value = '  +2.09%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.65'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.06%  '

$ws.Range("E11").Value = '  -0.18%  '

$ws.Range("E12").Value = '  -0.35%  '

$ws.Range("E13").Value = '  +0.49%  '

$ws.Range("D14").Value = '3.055.75'
$ws.Range("E14").Value = '  +1.57%  '

$ws.Range("D15").Value = '62.894.61'
$ws.Range("E15").Value = '  -0.11%  '

$ws.Range("E16").Value = '  +2.97%  '

$ws.Range("D17").Value = '2.599.16'
$ws.Range("E17").Value = '  +1.81%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.31'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.43%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '342.21'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.83%  '

$ws.Range("E20").Value = '  +1.28%  '

$ws.Range("E21").Value = '  -1.03%  '

$ws.Range("E22").Value = '  -0.09%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.24'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.15%  '

$ws.Range("D24").Value = '2.712.44'
$ws.Range("E24").Value = '  +1.41%  '

$ws.Range("B25").Value = 'Kaspa'
$ws.Range("C25").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.167'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.60%  '

$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.60'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.69%  '

$ws.Range("E27").Value = '  -0.28%  '

$ws.Range("E28").Value = '  -0.30%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.84'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +5.67%  '

$ws.Range("E30").Value = '  -1.50%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.93'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.94%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '471.40'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +14.93%  '

$ws.Range("E33").Value = '  +0.93%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '176.87'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.46%  '

$ws.Range("E35").Value = '  +4.30%  '

$ws.Range("E36").Value = '  +0.07%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.405'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.85%  '

$ws.Range("E38").Value = '  -0.52%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.52'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +4.00%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.70'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -2.45%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '158.39'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +4.85%  '

$ws.Range("E43").Value = '  -0.08%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.48'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +2.67%  '

$ws.Range("E45").Value = '  +5.14%  '

$ws.Range("E46").Value = '  +0.54%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0971'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.23%  '

$ws.Range("E48").Value = '  -1.23%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.36'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.30%  '

$ws.Range("E50").Value = '  +0.68%  '

$ws.Range("E51").Value = '  +1.05%  '
